$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, insert a new row 11 with a copy of what is currently row 10's data
# (this is the old row-10 record, shifted down to row 11 unchanged).
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(11, 3).Value = "Ñuble"
$ws.Cells.Item(11, 4).Value = 44832
$ws.Cells.Item(11, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11, 5).Value = 16
$ws.Cells.Item(11, 6).Value = 100114007
$ws.Cells.Item(11, 7).Value = "Jengibre"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 60
$ws.Cells.Item(11, 11).Value = 17000
$ws.Cells.Item(11, 12).Value = 18000
$ws.Cells.Item(11, 13).Value = 17500
$ws.Cells.Item(11, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(11, 15).Value = "Perú"
$ws.Cells.Item(11, 16).Value = 1346
$ws.Cells.Item(11, 17).Value = 13
$ws.Cells.Item(11, 18).Value = "Hortaliza"

# Update row 10 to hold what used to be row 9's data
$ws.Cells.Item(10, 4).Value = 44810
$ws.Cells.Item(10, 10).Value = 40
$ws.Cells.Item(10, 11).Value = 12000
$ws.Cells.Item(10, 12).Value = 13000
$ws.Cells.Item(10, 13).Value = 12500
$ws.Cells.Item(10, 16).Value = 962

# Update row 9 to hold the new record's data
$ws.Cells.Item(9, 4).Value = 44874
$ws.Cells.Item(9, 10).Value = 30
$ws.Cells.Item(9, 11).Value = 17000
$ws.Cells.Item(9, 12).Value = 17000
$ws.Cells.Item(9, 13).Value = 17000
$ws.Cells.Item(9, 16).Value = 1308
